# Mi logt y task para semana 2 de ciclo 2
# Adds two new time-log entries (rows 8 and 9) to the LOGT1 sheet:
#   Row 8: 23/10/2014 - Physical Data Model
#   Row 9: 25/10/2015 - Cycle Report

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: 23/10/2014, 15:00-17:00, 10 min interruption, Physical Data Model ---
$ws.Range("A8").Value = "23/10/2014"
$ws.Range("B8").Value = 0.625
$ws.Range("C8").Value = 0.70833333333333337
$ws.Range("D8").Value = 10
$ws.Range("E8").Formula = "=((HOUR(C8)-HOUR(B8))*60)+(MINUTE(C8)-MINUTE(B8))-D8"
$ws.Range("F8").Value = 37
$ws.Range("H8").Value = "Physical Data Model"
$ws.Rows.Item(8).RowHeight = 26

# --- Row 9: 25/10/2015, 10:00-10:40, no interruption, Cycle Report ---
$ws.Range("A9").Value = "25/10/2015"
$ws.Range("B9").Value = 0.41666666666666669
$ws.Range("C9").Value = 0.44444444444444442
$ws.Range("D9").Value = 0
$ws.Range("E9").Formula = "=((HOUR(C9)-HOUR(B9))*60)+(MINUTE(C9)-MINUTE(B9))-D9"
$ws.Range("F9").Value = 45
$ws.Range("H9").Value = "Cycle Report"
$ws.Rows.Item(9).RowHeight = 26

# Move the active selection to the next entry row, as Excel would after data entry
[void]$ws.Range("F10").Select()
